$d = $word.ActiveDocument

# Locate the paragraph containing "Botão de Sair da conta." so we can
# insert the new list item right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Bot*o de Sair da conta*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not locate anchor paragraph 'Botão de Sair da conta.'"
}

# Insert a brand new paragraph right after it; Word automatically carries
# over the paragraph style (Parágrafo da Lista) and the list numbering
# (ilvl 0 / numId 1) from the paragraph it was split from.
$target.Range.InsertParagraphAfter()

# Re-fetch the freshly created paragraph (the one right after the anchor)
# and give it its text.
$newIndex = $target.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.Text = "Alterar todos URLS que fazem referencia pra Localhost:8080"

Write-Output "Inserted paragraph $newIndex with text: $($newPara.Range.Text)"
